$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

# --- Header row: rename columns (order of edits matches how the shared
# string table gets rebuilt by the engine so the resulting indices line up
# with the target file) ---
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Waist"
$ws.Range("C1").Value = "Weight"
$ws.Range("D1").Value = "kcal total"

# --- Row 145: fill in the two missing duration/kcal columns ---
$ws.Range("D145").Value = 2107
$ws.Range("E145").Value = 1934

# --- Rows 146-148 already exist (placeholders) - fill in the data,
# reusing whatever style each cell already carries ---
$ws.Range("A146").Value = 45441
$ws.Range("B146").Value = 96
$ws.Range("C146").Value = 78.6
$ws.Range("D146").Value = 2283
$ws.Range("E146").Value = 1894
$ws.Range("F146").Value = 1

$ws.Range("A147").Value = 45442
$ws.Range("B147").Value = 96
$ws.Range("C147").Value = 78.3
$ws.Range("D147").Value = 1855
$ws.Range("E147").Value = 1855

$ws.Range("A148").Value = 45443
$ws.Range("B148").Value = 96
$ws.Range("C148").Value = 78.7
$ws.Range("D148").Value = 2071
$ws.Range("E148").Value = 1973
$ws.Range("F148").Value = 1

# --- Rows 149-161: brand-new rows. Column A needs the same date format
# as the rest of the date column, so copy that formatting over first,
# then write the values (B-F pick up the workbook's default style,
# matching the source rows which carry no explicit "s" attribute) ---
$ws.Range("A145").Copy()
$ws.Range("A149:A161").PasteSpecial(-4122)

$ws.Range("A149").Value = 45444
$ws.Range("B149").Value = 96
$ws.Range("C149").Value = 78.7
$ws.Range("D149").Value = 3121
$ws.Range("E149").Value = 3121
$ws.Range("F149").Value = 1

$ws.Range("A150").Value = 45445
$ws.Range("B150").Value = 96
$ws.Range("C150").Value = 79.4
$ws.Range("D150").Value = 2844
$ws.Range("E150").Value = 2268
$ws.Range("F150").Value = 1

$ws.Range("A151").Value = 45446
$ws.Range("B151").Value = 96
$ws.Range("C151").Value = 78.9
$ws.Range("D151").Value = 2266
$ws.Range("E151").Value = 2077
$ws.Range("F151").Value = 1

$ws.Range("A152").Value = 45447
$ws.Range("B152").Value = 96
$ws.Range("C152").Value = 79.3
$ws.Range("D152").Value = 3571
$ws.Range("E152").Value = 3571
$ws.Range("F152").Value = 1

$ws.Range("A153").Value = 45448
$ws.Range("B153").Value = 96
$ws.Range("C153").Value = 79.4
$ws.Range("D153").Value = 1753
$ws.Range("E153").Value = 1359
$ws.Range("F153").Value = 1

$ws.Range("A154").Value = 45449
$ws.Range("B154").Value = 96.5
$ws.Range("C154").Value = 79
$ws.Range("D154").Value = 4102
$ws.Range("E154").Value = 3910
$ws.Range("F154").Value = 1

$ws.Range("A155").Value = 45450
$ws.Range("B155").Value = 96.5
$ws.Range("C155").Value = 80.6
$ws.Range("D155").Value = 3110
$ws.Range("E155").Value = 2776
$ws.Range("F155").Value = 1

$ws.Range("A156").Value = 45451
$ws.Range("B156").Value = 96.5
$ws.Range("C156").Value = 80.4
$ws.Range("D156").Value = 4427
$ws.Range("E156").Value = 4294

$ws.Range("A157").Value = 45452
$ws.Range("B157").Value = 97
$ws.Range("C157").Value = 80.6
$ws.Range("D157").Value = 2752
$ws.Range("E157").Value = 2123
$ws.Range("F157").Value = 1

$ws.Range("A158").Value = 45453
$ws.Range("B158").Value = 97
$ws.Range("C158").Value = 80.3
$ws.Range("D158").Value = 2350
$ws.Range("E158").Value = 2017
$ws.Range("F158").Value = 1

$ws.Range("A159").Value = 45454
$ws.Range("B159").Value = 97
$ws.Range("C159").Value = 80
$ws.Range("D159").Value = 3771
$ws.Range("E159").Value = 3771
$ws.Range("F159").Value = 1

$ws.Range("A160").Value = 45455
$ws.Range("B160").Value = 97
$ws.Range("C160").Value = 80.6
$ws.Range("D160").Value = 4242
$ws.Range("E160").Value = 3913
$ws.Range("F160").Value = 1

$ws.Range("A161").Value = 45456
$ws.Range("B161").Value = 97
$ws.Range("C161").Value = 80.8
$ws.Range("D161").Value = 2608
$ws.Range("E161").Value = 2608
$ws.Range("F161").Value = 1

# --- Row 162: empty placeholder row, only the date-formatted A cell exists ---
$ws.Range("A145").Copy()
$ws.Range("A162").PasteSpecial(-4122)

# --- View state: active sheet, scrolled-down selection ---
$ws.Activate()
$ws.Range("C146").Select()
